# Work on Tef Model
# - Rename CWP2011PopHighTreatment1..16 (A44:A59) -> CWP2011Treatment1..16
# - Re-scale observed harvest data (Biomass/GrainWt/StrawWt) from kg/ha-ish
#   units down by a factor of 10 (D, E columns), and replace the StrawWt
#   (F column) formula "=D-E" with its literal computed value, also scaled.
# - Update the saved view state (pane/selection) for the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ObservedHarvestData")

# --- 1. Column A relabeling for rows 44-59 ---------------------------------
for ($i = 1; $i -le 16; $i++) {
    $row = 43 + $i
    $ws.Cells.Item($row, 1).Value = "CWP2011Treatment$i"
}

# --- 2. Rescale D (Biomass), E (GrainWt) and de-formularize F (StrawWt) ---
for ($r = 2; $r -le 83; $r++) {
    $d = $ws.Cells.Item($r, 4).Value()
    $e = $ws.Cells.Item($r, 5).Value()
    $f = $ws.Cells.Item($r, 6).Value()

    if ($f -ne $null) {
        $ws.Cells.Item($r, 6).Value = $f / 10
    }
    if ($d -ne $null) {
        $ws.Cells.Item($r, 4).Value = $d / 10
    }
    if ($e -ne $null) {
        $ws.Cells.Item($r, 5).Value = $e / 10
    }
}

# --- 3. View state: scroll/selection -------------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 37
$win.ScrollColumn = 2
$ws.Range("Q69").Select()
